$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S3").Value = 1.3

$ws.Range("G4").Value = 1.4
$ws.Range("I4").Value = 7.5
$ws.Range("J4").Value = 1.95
$ws.Range("L4").Value = 7.5
$ws.Range("Q4").Value = 1.99
$ws.Range("R4").Value = 1.91
$ws.Range("S4").Value = 1.37
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 9
$ws.Range("AD4").Value = 8.5
$ws.Range("AG4").Value = 17
$ws.Range("AM4").Value = 501
$ws.Range("AO4").Value = 7
$ws.Range("AW4").Value = 8.5

$ws.Range("S5").Value = 1.37

$ws.Range("G6").Value = 2.55
$ws.Range("J6").Value = 3.25
$ws.Range("S6").Value = 1.5
$ws.Range("T6").Value = 2.37
$ws.Range("Y6").Value = 10
$ws.Range("AD6").Value = 6
$ws.Range("AJ6").Value = 34
$ws.Range("AN6").Value = 4.33
$ws.Range("AX6").Value = 19
$ws.Range("AZ6").Value = 67

$ws.Range("O9").Value = 1.67
$ws.Range("P9").Value = 2.1
